$wb = $excel.ActiveWorkbook

# Rename the "Include from Health Data Conn" sheet to "Include from Medication Knowl"
$includeSheet = $wb.Worksheets.Item("Include from Health Data Conn")
$includeSheet.Name = "Include from Medication Knowl"

# Update the Metadata sheet values
$metaSheet = $wb.Worksheets.Item("Metadata")

$metaSheet.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/ValueSet/medication-knowledge-cost-type"
$metaSheet.Range("B3").Value = "8.0.0"
$metaSheet.Range("B4").Value = "MedicationKnowledgeCostTypValueSet"
$metaSheet.Range("B5").Value = "Medication Knowledge Cost Type Value Set"
$metaSheet.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$metaSheet.Range("B9").Value = "LinuxForHealth Team"
$metaSheet.Range("B11").Value = "Medication knowledge cost type value set"

# Update the System URI value on the Include sheet
$includeSheet.Range("B4").Value = "http://linuxforhealth.org/fhir/cdm/CodeSystem/medication-knowledge-cost-type"
